$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in the header (A1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 10:05"

# Refresh per-country statistics (new snapshot) and restore descending sort
# order by "Casos totales" (column B) for the rows whose rank changed.

# Row 6: Rusia
$ws.Range("A6").Value = "Rusia"
$ws.Range("B6").Value = 405843
$ws.Range("C6").Value = 9268
$ws.Range("D6").Value = 171883
$ws.Range("E6").Value = 229267
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 138
$ws.Range("H6").Value = 4693

# Row 11: Alemania
$ws.Range("A11").Value = "Alemania"
$ws.Range("B11").Value = 183294
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 165200
$ws.Range("E11").Value = 9494
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 8600

# Row 29: Singapur
$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 34884
$ws.Range("C29").Value = 518
$ws.Range("D29").Value = 20727
$ws.Range("E29").Value = 14134
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 23

# Row 38: Polonia
$ws.Range("A38").Value = "Polonia"
$ws.Range("B38").Value = 23686
$ws.Range("C38").Value = 115
$ws.Range("D38").Value = 11271
$ws.Range("E38").Value = 11351
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 1064

# Row 39: Ucrania
$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 23672
$ws.Range("C39").Value = 468
$ws.Range("D39").Value = 9538
$ws.Range("E39").Value = 13426
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 12
$ws.Range("H39").Value = 708

# Row 41: Rumania
$ws.Range("A41").Value = "Rumania"
$ws.Range("B41").Value = 19133
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 13046
$ws.Range("E41").Value = 4825
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 1262

# Row 55: Oman
$ws.Range("A55").Value = "Oman"
$ws.Range("B55").Value = 10423
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 2396
$ws.Range("E55").Value = 7983
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 44

# Row 99: Lituania
$ws.Range("A99").Value = "Lituania"
$ws.Range("B99").Value = 1675
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 1236
$ws.Range("E99").Value = 369
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 70

# Row 100: Maldivas
$ws.Range("A100").Value = "Maldivas"
$ws.Range("B100").Value = 1672
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 406
$ws.Range("E100").Value = 1261
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 5

# Row 126: Georgia
$ws.Range("A126").Value = "Georgia"
$ws.Range("B126").Value = 783
$ws.Range("C126").Value = 26
$ws.Range("D126").Value = 605
$ws.Range("E126").Value = 166
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 12

# Row 127: Principado de Andorra
$ws.Range("A127").Value = "Principado de Andorra"
$ws.Range("B127").Value = 764
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 692
$ws.Range("E127").Value = 21
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 51

# Row 128: Nicaragua
$ws.Range("A128").Value = "Nicaragua"
$ws.Range("B128").Value = 759
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 370
$ws.Range("E128").Value = 354
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 35

# Row 129: Republica del Chad
$ws.Range("A129").Value = "Republica del Chad"
$ws.Range("B129").Value = 759
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 470
$ws.Range("E129").Value = 224
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 65

# Row 130: Madagascar
$ws.Range("A130").Value = "Madagascar"
$ws.Range("B130").Value = 758
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 165
$ws.Range("E130").Value = 587
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 6

# Row 200: Belice
$ws.Range("A200").Value = "Belice"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 16
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 2

# Row 201: Santa Lucia
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 18
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# Row 213: Islas Virgenes Britanicas
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Papua Nueva Guinea
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
